# Generate Report for Handoff
# Replace the old GUID-based file name ("cdd60225-56a7-4041-a34a-0b46b127ed22")
# with the new one ("acc6e68f-631e-43f6-a6c7-ea97c063043d") and the new xliff
# hash ("69fdb3b38349d331b5c4f961d7353d9277d17f3a" -> "2f8f1c07ecb678d0e488e04d22bd1c1d9a7fe724")
# across the Overview / zh-cn / de-de sheets, and bump the recorded timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "cdd60225-56a7-4041-a34a-0b46b127ed22"
$newGuid = "acc6e68f-631e-43f6-a6c7-ea97c063043d"
$oldHash = "69fdb3b38349d331b5c4f961d7353d9277d17f3a"
$newHash = "2f8f1c07ecb678d0e488e04d22bd1c1d9a7fe724"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-20 14:59:45"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-20 14:59:41"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-20 14:59:45"
